$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row (row 11): Right marks 5 -> 4, Wrong marks -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Update "Total" row (row 12): Right total 120 -> 96, Wrong total -1 -> -2
$ws.Range("B12").Value = 96
$ws.Range("C12").Value = -2

# Update the Max column display text "120 / 140" -> "94 / 112"
$ws.Range("E12").Value = "94 / 112"
